$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 41
# from 45224 to 45233 (serial date values).
$ws.Range("C2:C41").Value = 45233
